$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire first row (Manchester United v Liverpool / 24 OctSun16:30).
# All rows below shift up by one, and the sheet's used range shrinks from
# A1:B16 to A1:B15.
$ws.Rows("1:1").Delete()
